# TC for Order tracking and cancelltion updated
#
# Updates the " Order Tracking & Cancellation" sheet (sheet10.xml):
#  - shifts the existing tracking/cancellation test descriptions down one row
#    to make room for a new "Verify order tracking" scenario
#  - adds two new rows (TS_TRACK_05 / TS_TRACK_06) for the scenarios that
#    were displaced
#  - makes this sheet the active/selected sheet (it was previously the
#    "Order Confirmation" sheet that was active)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(" Order Tracking & Cancellation")

# Row 8: new "Verify order tracking" scenario description (TS_TRACK_01 / FR-TRACK-01 stay the same)
$ws.Cells.Item(8, 3).Value = "Verify order tracking "

# Row 9: shift "Verify order tracking status" down from row 8 (TS_TRACK_02 / FR-TRACK-02 stay the same)
$ws.Cells.Item(9, 3).Value = "Verify order tracking status"

# Row 10: shift "Verify order cancellation before shipment" down from row 9 (TS_TRACK_03 / FR-TRACK-03 stay the same)
$ws.Cells.Item(10, 3).Value = "Verify order cancellation before shipment"

# Row 11: new description replacing "Verify refund initiation after cancellation" (TS_TRACK_04 / FR-TRACK-04 stay the same)
$ws.Cells.Item(11, 3).Value = "Verify  to Prevent cancellation after shipment"

# New requirement IDs for the two displaced scenarios (rows 12 & 13)
$ws.Cells.Item(12, 2).Value = "FR-TRACK-05"
$ws.Cells.Item(13, 2).Value = "FR-TRACK-06"

# New scenario IDs for the two displaced scenarios (rows 12 & 13)
$ws.Cells.Item(12, 1).Value = "TS_TRACK_05"
$ws.Cells.Item(13, 1).Value = "TS_TRACK_06"

# Row 12 (new): displaced "Verify cancellation confirmation message" scenario
$ws.Cells.Item(12, 3).Value = "Verify cancellation confirmation message"

# Row 13 (new): displaced "Verify refund initiation after cancellation" scenario
$ws.Cells.Item(13, 3).Value = "Verify refund initiation after cancellation"

# Make this the active sheet / selected cell, matching the saved view state
$ws.Activate()
$ws.Range("C15").Select()
